# The deck's Slide Master currently carries the "Integral" design with the
# "Red Violet" colour scheme (ppt/theme/theme2.xml). The edit restores the
# default Office colour scheme (originally ppt/theme/theme1.xml, which is
# only linked from the Notes Master and isn't reachable through the
# PowerPoint object model) on the Slide Master's theme.
#
# PowerPoint exposes the twelve theme colour slots - dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink, in that fixed order - through
# Slide.ThemeColorScheme; writing .RGB there edits the OOXML theme part
# used by the active Slide Master. RGB values use the usual VBA encoding:
# R + G*256 + B*65536.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
